# Update cryptocurrency price/volume figures (columns D and E) to match
# the latest GitHub Actions scrape. Values are written as text (not
# numbers) so formats like "29.970.71" / "0.3950" survive untouched,
# mirroring the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    # Force text interpretation so numeric-looking strings (e.g. "42.21")
    # are not coerced into numbers, then drop back to the default style
    # so no stray formatting is left behind on the cell.
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '29.970.71'
Set-TextValue 'E2' '  -0.36%  '
Set-TextValue 'D3' '1.870.97'
Set-TextValue 'E3' '  -2.62%  '
Set-TextValue 'E4' '  +0.08%  '
Set-TextValue 'D5' '319.23'
Set-TextValue 'E5' '  -3.82%  '
Set-TextValue 'E7' '  -2.81%  '
Set-TextValue 'D8' '0.3950'
Set-TextValue 'E8' '  -2.51%  '
Set-TextValue 'D9' '0.08207'
Set-TextValue 'E9' '  -3.66%  '
Set-TextValue 'D10' '42.21'
Set-TextValue 'E10' '  -1.82%  '
Set-TextValue 'D11' '1.094'
Set-TextValue 'E11' '  -3.05%  '
Set-TextValue 'D12' '22.91'
Set-TextValue 'E12' '  +2.44%  '
Set-TextValue 'D13' '1.863.36'
Set-TextValue 'E13' '  -2.92%  '
Set-TextValue 'E14' '  -1.86%  '
Set-TextValue 'D15' '7.201'
Set-TextValue 'E15' '  -2.76%  '
Set-TextValue 'D16' '1.002'
Set-TextValue 'E16' '  +0.16%  '
Set-TextValue 'D17' '91.92'
Set-TextValue 'E17' '  -4.74%  '
Set-TextValue 'E18' '  -2.71%  '
Set-TextValue 'D19' '0.06372'
Set-TextValue 'E19' '  -4.99%  '
Set-TextValue 'E20' '  -2.03%  '
Set-TextValue 'E21' '  +0.04%  '
Set-TextValue 'D22' '29.974.01'
Set-TextValue 'E22' '  -0.36%  '
Set-TextValue 'D23' '5.830'
Set-TextValue 'E23' '  -3.91%  '
Set-TextValue 'D24' '11.12'
Set-TextValue 'E24' '  -1.16%  '
Set-TextValue 'D25' '2.172'
Set-TextValue 'E25' '  -2.49%  '
Set-TextValue 'D26' '2.088.75'
Set-TextValue 'E26' '  -2.40%  '
Set-TextValue 'D27' '161.03'
Set-TextValue 'E27' '  +0.55%  '
Set-TextValue 'D28' '20.95'
Set-TextValue 'E28' '  -0.90%  '
Set-TextValue 'D29' '2.240'
Set-TextValue 'E29' '  -8.77%  '
Set-TextValue 'D30' '127.41'
Set-TextValue 'E30' '  -1.75%  '
Set-TextValue 'D31' '1.070'
Set-TextValue 'E31' '  -1.24%  '
Set-TextValue 'E32' '  -2.37%  '
Set-TextValue 'D33' '5.934'
Set-TextValue 'E33' '  -2.93%  '
Set-TextValue 'E34' '  +2.14%  '
Set-TextValue 'E35' '  -3.70%  '
Set-TextValue 'D36' '5.217'
Set-TextValue 'E36' '  -0.10%  '
Set-TextValue 'D37' '0.06370'
Set-TextValue 'E37' '  -3.74%  '
Set-TextValue 'E38' '  -4.02%  '
Set-TextValue 'E39' '  -5.23%  '
Set-TextValue 'E40' '  -5.81%  '
Set-TextValue 'D41' '0.6304'
Set-TextValue 'E41' '  -4.09%  '
Set-TextValue 'E42' '  -3.66%  '
Set-TextValue 'D43' '1.207'
Set-TextValue 'E43' '  -3.10%  '
Set-TextValue 'D44' '0.9998'
Set-TextValue 'E44' '  +0.06%  '
Set-TextValue 'D45' '0.5910'
Set-TextValue 'E45' '  -4.86%  '
Set-TextValue 'D46' '12.90'
Set-TextValue 'E46' '  -2.94%  '
Set-TextValue 'D47' '3.634'
Set-TextValue 'E47' '  -4.17%  '
Set-TextValue 'D48' '2.005'
Set-TextValue 'E48' '  -3.92%  '
Set-TextValue 'D49' '122.59'
Set-TextValue 'E49' '  -2.46%  '
Set-TextValue 'E51' '  -2.94%  '
